$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: handback status messages + widened status columns (E, F)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet: refresh handback datetime, clear stale error detail
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("K2").Value = "2016-10-20 00:14:15"
$wsZh.Range("K3").Value = "2016-10-20 00:14:15"
$wsZh.Range("P2").Value = ""
$wsZh.Range("P3").Value = ""
$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------------
# de-de sheet: generated handback xliff, refreshed handback datetime, clear
# stale error detail
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-20 00:14:33"
$wsDe.Range("P2").Value = ""
$wsDe.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-10-20 00:14:33"
$wsDe.Range("P3").Value = ""
$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(16).ColumnWidth = 12.833333333333334
